$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "prepare" benchmarking pass produced 12 more result rows that were
# missing from the sheet (row 463-474). Append them below the existing
# data. Copy row 462 down first so the new rows inherit its number
# formats/styles (columns Q:W -> 0.000000, X -> 0.000, Y:Z -> 0.00), then
# stamp every cell with the recorded values.
$ws.Range("A462:Z462").Copy($ws.Range("A463:Z474"))

$newRows = @(
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 2.286108016967773, 0.009144432067871093, 0.002155542373657227, 0.003625, 0.0018125, 0, 35.2994372844696, 4268.8, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 2.393219947814941, 0.009572879791259765, 0.001131772994995117, 0.004757812500000001, 0.00237890625, 0, 37.41122794151306, 4756, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 2.524171352386475, 0.0100966854095459, 0.001620769500732422, 0.004757812500000001, 0.00237890625, 0, 37.41674160957336, 4756, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 2.611611366271973, 0.01044644546508789, 0.00184941291809082, 0.0047578125, 0.00237890625, 0, 37.41062450408936, 4756, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 2.542445659637451, 0.01016978263854981, 0.001868247985839844, 0.004531250000000001, 0.002265625, 0, 38.1127336025238, 4918.4, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 2.240594387054443, 0.008962377548217772, 0.001659870147705078, 0.0038515625, 0.00192578125, 0, 37.00537061691284, 4663.2, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 2.207418918609619, 0.008829675674438476, 0.001140594482421875, 0.004078125, 0.0020390625, 0, 37.91260147094727, 4872, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 2.156326770782471, 0.008625307083129884, 0.0004904270172119141, 0.004531250000000001, 0.002265625, 0, 35.68000936508179, 4361.6, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 1.548069953918457, 0.006192279815673829, 0.0004827976226806641, 0.0029453125, 0.00147265625, 0, 34.9842791557312, 4199.2, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 1.712746620178223, 0.006850986480712891, 0.0004765987396240234, 0.0033984375, 0.00169921875, 0, 36.80687069892883, 4616.8, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 2.321778297424316, 0.009287113189697266, 0.0004842281341552734, 0.004984375, 0.0024921875, 0, 38.35043406486511, 4964.8, 928),
    @(2, 1, 2, 10, 10, 3, 50, 0, 0, 1000, 2000, $false, 16, 16, '<-parameter / result->', 4, -1, 1.778490543365479, 0.007113962173461915, 0.0006968975067138672, 0.0033984375, 0.00169921875, 0, 36.82202172279358, 4616.8, 928)
)

$startRow = 463
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
